# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Moogle_Profits workbook
# (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1429.3334
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 1358
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 1358
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -2010

# Row 40
$ws.Range("H40").Value = 3957.0417
$ws.Range("J40").Value = 4467.625
$ws.Range("L40").Value = 4467.625
$ws.Range("N40").Value = -4817.625

# Row 43
$ws.Range("H43").Value = 3914.8462
$ws.Range("I43").Value = 2600
$ws.Range("J43").Value = 4499.222
$ws.Range("K43").Value = 2600
$ws.Range("L43").Value = 4499.222
$ws.Range("M43").Value = -2531
$ws.Range("N43").Value = -4637.222

# Row 51
$ws.Range("H51").Value = 7449.5
$ws.Range("J51").Value = 7449.5
$ws.Range("L51").Value = 7449.5
$ws.Range("N51").Value = -8417.5

# Row 70
$ws.Range("H70").Value = 2706.375
$ws.Range("J70").Value = 3730.2
$ws.Range("L70").Value = 11190.6
$ws.Range("N70").Value = -11730.6

# Row 73
$ws.Range("H73").Value = 2706.375
$ws.Range("J73").Value = 3730.2
$ws.Range("L73").Value = 11190.6
$ws.Range("N73").Value = -13062.6

# Row 80
$ws.Range("H80").Value = 1060.875
$ws.Range("I80").Value = 848.5454999999999
$ws.Range("K80").Value = 2545.6365
$ws.Range("M80").Value = -1547.6365

# Row 83
$ws.Range("H83").Value = 1060.875
$ws.Range("I83").Value = 848.5454999999999
$ws.Range("K83").Value = 7636.9095
$ws.Range("M83").Value = -2644.9095

# Row 113
$ws.Range("H113").Value = 4999.3335
$ws.Range("J113").Value = 4999.3335
$ws.Range("L113").Value = 4999.3335
$ws.Range("N113").Value = -11507.3335

# Row 138
$ws.Range("H138").Value = 4330.1567
$ws.Range("I138").Value = 2429.7083
$ws.Range("J138").Value = 6019.4443
$ws.Range("K138").Value = 7289.124899999999
$ws.Range("L138").Value = 18058.3329
$ws.Range("M138").Value = -2149.124899999999
$ws.Range("N138").Value = -28338.3329

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6709.3076
$ws.Range("I32").Value = 4518.5835
$ws.Range("K32").Value = 4518.5835
$ws.Range("M32").Value = -4231.5835

# Row 74
$ws.Range("H74").Value = 4609.6665
$ws.Range("I74").Value = 2528.5
$ws.Range("K74").Value = 2528.5
$ws.Range("M74").Value = -1654.5

# Row 77
$ws.Range("H77").Value = 4609.6665
$ws.Range("I77").Value = 2528.5
$ws.Range("K77").Value = 12642.5
$ws.Range("M77").Value = -8274.5

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 6171.2354
$ws.Range("I86").Value = 3638.6667
$ws.Range("J86").Value = 12249.4
$ws.Range("K86").Value = 3638.6667
$ws.Range("L86").Value = 12249.4
$ws.Range("M86").Value = -2515.6667
$ws.Range("N86").Value = -14495.4

# Row 89
$ws.Range("H89").Value = 6171.2354
$ws.Range("I89").Value = 3638.6667
$ws.Range("J89").Value = 12249.4
$ws.Range("K89").Value = 18193.3335
$ws.Range("L89").Value = 61247
$ws.Range("M89").Value = -12577.3335
$ws.Range("N89").Value = -72479

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5280.967
$ws.Range("I31").Value = 3210.075
$ws.Range("K31").Value = 3210.075
$ws.Range("M31").Value = -2915.075

# Row 34
$ws.Range("H34").Value = 5280.967
$ws.Range("I34").Value = 3210.075
$ws.Range("K34").Value = 3210.075
$ws.Range("M34").Value = -3008.075

# Row 62
$ws.Range("H62").Value = 16674446
$ws.Range("I62").Value = 7619.7
$ws.Range("J62").Value = 50008100
$ws.Range("K62").Value = 7619.7
$ws.Range("L62").Value = 50008100
$ws.Range("M62").Value = -6995.7
$ws.Range("N62").Value = -50009348

# Row 65
$ws.Range("H65").Value = 16674446
$ws.Range("I65").Value = 7619.7
$ws.Range("J65").Value = 50008100
$ws.Range("K65").Value = 38098.5
$ws.Range("L65").Value = 250040500
$ws.Range("M65").Value = -34978.5
$ws.Range("N65").Value = -250046740

# Row 81
$ws.Range("H81").Value = 199950
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

# Row 84
$ws.Range("H84").Value = 199950
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

# Row 134
$ws.Range("H134").Value = 5113.222
$ws.Range("I134").Value = 3502.375
$ws.Range("J134").Value = 18000
$ws.Range("K134").Value = 10507.125
$ws.Range("L134").Value = 54000
$ws.Range("M134").Value = -7972.125
$ws.Range("N134").Value = -59070

$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Range("H63").Value = 18724.25
$ws.Range("I63").Value = 21632.666
$ws.Range("K63").Value = 64897.99800000001
$ws.Range("M63").Value = -64148.99800000001

# Row 66
$ws.Range("H66").Value = 18724.25
$ws.Range("I66").Value = 21632.666
$ws.Range("K66").Value = 194693.994
$ws.Range("M66").Value = -190949.994

# Row 112
$ws.Range("H112").Value = 11993.333
$ws.Range("I112").Value = 8441.333000000001
$ws.Range("K112").Value = 25323.999
$ws.Range("M112").Value = -24215.999

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 478.8
$ws.Range("I107").Value = 398.33334
$ws.Range("K107").Value = 398.33334
$ws.Range("M107").Value = 1521.66666

# Row 113
$ws.Range("H113").Value = 2447.3823
$ws.Range("J113").Value = 3354.3845
$ws.Range("L113").Value = 3354.3845
$ws.Range("N113").Value = -7694.3845

# Row 136
$ws.Range("H136").Value = 25032.25
$ws.Range("J136").Value = 25032.25
$ws.Range("L136").Value = 75096.75
$ws.Range("N136").Value = -80196.75

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 5785.9443
$ws.Range("I68").Value = 5743.1333
$ws.Range("K68").Value = 5743.1333
$ws.Range("M68").Value = -4994.1333

# Row 71
$ws.Range("H71").Value = 5785.9443
$ws.Range("I71").Value = 5743.1333
$ws.Range("K71").Value = 28715.6665
$ws.Range("M71").Value = -24971.6665

# Row 96
$ws.Range("H96").Value = 139093.12
$ws.Range("J96").Value = 139093.12
$ws.Range("L96").Value = 139093.12
$ws.Range("N96").Value = -144585.12

# Row 98
$ws.Range("H98").Value = 161211.25
$ws.Range("J98").Value = 161211.25
$ws.Range("L98").Value = 161211.25
$ws.Range("N98").Value = -167201.25

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1703.5294
$ws.Range("I126").Value = 1642.6428
$ws.Range("K126").Value = 4927.928400000001
$ws.Range("M126").Value = -2457.928400000001

# Row 130
$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040

